# Applies market-price / profit updates scraped by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 4749.625
$ws.Range("I46").Value = 2999.25
$ws.Range("J46").Value = 6500
$ws.Range("K46").Value = 8997.75
$ws.Range("L46").Value = 19500
$ws.Range("M46").Value = -8878.75
$ws.Range("N46").Value = -19738

$ws.Range("H60").Value = 4749.625
$ws.Range("I60").Value = 2999.25
$ws.Range("J60").Value = 6500
$ws.Range("K60").Value = 8997.75
$ws.Range("L60").Value = 19500
$ws.Range("M60").Value = -8513.75
$ws.Range("N60").Value = -20468

$ws.Range("H86").Value = 4639.851
$ws.Range("I86").Value = 3971.1516
$ws.Range("K86").Value = 3971.1516
$ws.Range("M86").Value = -2848.1516

$ws.Range("H88").Value = 3083.9285
$ws.Range("J88").Value = 2134.1667
$ws.Range("L88").Value = 2134.1667
$ws.Range("N88").Value = -2946.1667

$ws.Range("H89").Value = 4639.851
$ws.Range("I89").Value = 3971.1516
$ws.Range("K89").Value = 19855.758
$ws.Range("M89").Value = -14239.758

$ws.Range("H91").Value = 3083.9285
$ws.Range("J91").Value = 2134.1667
$ws.Range("L91").Value = 2134.1667
$ws.Range("N91").Value = -4942.1667

$ws.Range("H106").Value = 2976.4443
$ws.Range("I106").Value = 3568.5715
$ws.Range("K106").Value = 3568.5715
$ws.Range("M106").Value = -2937.5715

$ws.Range("H132").Value = 29749.572
$ws.Range("I132").Value = 29749.572
$ws.Range("K132").Value = 89248.716
$ws.Range("M132").Value = -86718.716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1143.6
$ws.Range("I2").Value = 948.4
$ws.Range("J2").Value = 1338.8
$ws.Range("K2").Value = 948.4
$ws.Range("L2").Value = 1338.8
$ws.Range("M2").Value = -835.4
$ws.Range("N2").Value = -1564.8

$ws.Range("H32").Value = 1295870.5
$ws.Range("I32").Value = 1463542.5
$ws.Range("K32").Value = 1463542.5
$ws.Range("M32").Value = -1463255.5

$ws.Range("H57").Value = 4000
$ws.Range("I57").Value = 4000
$ws.Range("K57").Value = 4000
$ws.Range("M57").Value = -3516

$ws.Range("H61").Value = 6253154.5
$ws.Range("I61").Value = 2547.9
$ws.Range("K61").Value = 2547.9
$ws.Range("M61").Value = -2335.9

$ws.Range("H92").Value = 100000
$ws.Range("J92").Value = 100000
$ws.Range("L92").Value = 100000
$ws.Range("N92").Value = -104992

$ws.Range("H116").Value = 1143.6
$ws.Range("I116").Value = 948.4
$ws.Range("J116").Value = 1338.8
$ws.Range("K116").Value = 948.4
$ws.Range("L116").Value = 1338.8
$ws.Range("M116").Value = 1345.6
$ws.Range("N116").Value = -5926.8

$ws.Range("H122").Value = 1896.25
$ws.Range("I122").Value = 1773.4445
$ws.Range("J122").Value = 2264.6667
$ws.Range("K122").Value = 5320.333500000001
$ws.Range("L122").Value = 6794.000100000001
$ws.Range("M122").Value = -2870.333500000001
$ws.Range("N122").Value = -11694.0001

$ws.Range("H132").Value = 3096.2222
$ws.Range("I132").Value = 1455
$ws.Range("J132").Value = 5483.4546
$ws.Range("K132").Value = 4365
$ws.Range("L132").Value = 16450.3638
$ws.Range("M132").Value = -1835
$ws.Range("N132").Value = -21510.3638

$ws.Range("H134").Value = 60000.668
$ws.Range("J134").Value = 60000.668
$ws.Range("L134").Value = 60000.668
$ws.Range("N134").Value = -70140.66800000001

$ws.Range("H136").Value = 6253154.5
$ws.Range("I136").Value = 2547.9
$ws.Range("K136").Value = 7643.700000000001
$ws.Range("M136").Value = -5093.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1143.6
$ws.Range("I3").Value = 948.4
$ws.Range("J3").Value = 1338.8
$ws.Range("K3").Value = 948.4
$ws.Range("L3").Value = 1338.8
$ws.Range("M3").Value = -834.4
$ws.Range("N3").Value = -1566.8

$ws.Range("H86").Value = 3225.5
$ws.Range("I86").Value = 2114.1428
$ws.Range("K86").Value = 2114.1428
$ws.Range("M86").Value = -991.1428000000001

$ws.Range("H89").Value = 3225.5
$ws.Range("I89").Value = 2114.1428
$ws.Range("K89").Value = 10570.714
$ws.Range("M89").Value = -4954.714

$ws.Range("H105").Value = 5743.3335
$ws.Range("I105").Value = 2255.8572
$ws.Range("J105").Value = 17949.5
$ws.Range("K105").Value = 2255.8572
$ws.Range("L105").Value = 17949.5
$ws.Range("M105").Value = -508.8571999999999
$ws.Range("N105").Value = -21443.5

$ws.Range("H134").Value = 2689592.8
$ws.Range("I134").Value = 1397.931
$ws.Range("K134").Value = 4193.793
$ws.Range("M134").Value = -1658.793

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 169253.5
$ws.Range("I16").Value = 2877
$ws.Range("J16").Value = 502006.5
$ws.Range("K16").Value = 2877
$ws.Range("L16").Value = 502006.5
$ws.Range("M16").Value = -2590
$ws.Range("N16").Value = -502580.5

$ws.Range("H113").Value = 169253.5
$ws.Range("I113").Value = 2877
$ws.Range("J113").Value = 502006.5
$ws.Range("K113").Value = 2877
$ws.Range("L113").Value = 502006.5
$ws.Range("M113").Value = -707
$ws.Range("N113").Value = -506346.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 109.71429
$ws.Range("I12").Value = 18
$ws.Range("J12").Value = 134.72728
$ws.Range("K12").Value = 54
$ws.Range("L12").Value = 404.18184
$ws.Range("M12").Value = 119
$ws.Range("N12").Value = -750.18184

$ws.Range("H121").Value = 12445.348
$ws.Range("J121").Value = 17621.125
$ws.Range("L121").Value = 52863.375
$ws.Range("N121").Value = -55483.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 75.16
$ws.Range("I2").Value = 65
$ws.Range("J2").Value = 115.8
$ws.Range("K2").Value = 65
$ws.Range("L2").Value = 115.8
$ws.Range("M2").Value = 48
$ws.Range("N2").Value = -341.8

$ws.Range("H102").Value = 1708.5
$ws.Range("I102").Value = 1766.875
$ws.Range("J102").Value = 1475
$ws.Range("K102").Value = 1766.875
$ws.Range("L102").Value = 1475
$ws.Range("M102").Value = -144.875
$ws.Range("N102").Value = -4719

$ws.Range("H132").Value = 14212.071
$ws.Range("I132").Value = 7766.9565
$ws.Range("K132").Value = 23300.8695
$ws.Range("M132").Value = -20770.8695

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1436.3235
$ws.Range("J55").Value = 1495.421
$ws.Range("L55").Value = 1495.421
$ws.Range("N55").Value = -1841.421

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H130").Value = 58499.75
$ws.Range("J130").Value = 58499.75
$ws.Range("L130").Value = 58499.75
$ws.Range("N130").Value = -68539.75

$ws.Range("H132").Value = 1452763.9
$ws.Range("I132").Value = 3033733.5
$ws.Range("J132").Value = 3541.75
$ws.Range("K132").Value = 9101200.5
$ws.Range("L132").Value = 10625.25
$ws.Range("M132").Value = -9098670.5
$ws.Range("N132").Value = -15685.25

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6947176.5
$ws.Range("I132").Value = 7578377
$ws.Range("K132").Value = 22735131
$ws.Range("M132").Value = -22732601
